$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @{A=44820; B="KA53MD2318"; C="MICRA"; D="PMS"; E="WORK DONE DELIVERED"; F=6203; G=$null},
    @{A=44821; B="KA03MP9766"; C="FIGO"; D="RUNNING REPAIR"; E="WORK DONE DELIVERED"; F=2860; G="CREDIT"},
    @{A=44821; B="KA03MV0746"; C="ERTIGA"; D="PMS"; E="WORK DONE DELIVERED"; F=10823; G="CREDIT"},
    @{A=44821; B="KA03MU4168"; C="CELERIO"; D="PMS"; E="WORK DONE DELIVERED"; F=3403; G="CREDIT"},
    @{A=44821; B="KA03MU1016"; C="SWIFT"; D="RUNNING REPAIR"; E="WORK DONE DELIVERED"; F=12908; G=$null},
    @{A=44821; B="AP29AT0483"; C="RITZ"; D="RUNNING REPAIR"; E="WORK DONE DELIVERED"; F=2926; G=$null},
    @{A=44821; B="KA03NE7365"; C="NEXON"; D="PMS"; E="WORK DONE DELIVERED"; F=4378; G=$null},
    @{A=44823; B="KA03MM7229"; C="I10"; D="RUNNING REPAIR"; E="WORK DONE DELIVERED"; F=5553; G="CREDIT"},
    @{A=44823; B="KA03MX9007"; C="POLO"; D="BODY SHOP"; E="WORK DONE DELIVERED"; F=15291; G="CARD"},
    @{A=44823; B="TS07GQ1187"; C="WR-V"; D="PMS"; E="WORK DONE DELIVERED"; F=7526; G=$null},
    @{A=44823; B="KA01MR3476"; C="TIAGO"; D="PMS"; E="WORK DONE DELIVERED"; F=4930; G=$null},
    @{A=44824; B="KA01MK1436"; C="POLO"; D="PMS"; E="WORK DONE DELIVERED"; F=28053; G=$null},
    @{A=44824; B="KA03ME6233"; C="CRETA"; D="PMS"; E="WORK DONE DELIVERED"; F=6016; G=$null},
    @{A=44824; B="KA03MP9135"; C="FIGO"; D="BODY SHOP"; E="WORK DONE DELIVERED"; F=12748; G=$null},
    @{A=44825; B="KA53MA5526"; C="DUSTER"; D="PMS"; E="WORK DONE DELIVERED"; F=4424; G=$null},
    @{A=44825; B="DL8CAM5454"; C="SWIFT"; D="PMS"; E="WORK DONE DELIVERED"; F=4912; G=$null},
    @{A=44825; B="KA51MA9028"; C="FIGO"; D="PMS"; E="WORK DONE DELIVERED"; F=22130; G=$null},
    @{A=44825; B="TN57BA3434"; C="I20"; D="PMS"; E="WORK DONE DELIVERED"; F=13611; G=$null},
    @{A=44826; B="DL83AE7248"; C="XCENT"; D="PMS"; E="WORK DONE DELIVERED"; F=10599; G="CREDIT"},
    @{A=44826; B="KA03MS1975"; C="ETIOS"; D="RUNNING REPAIRR"; E="WORK DONE DELIVERED"; F=12711; G="CREDIT"},
    @{A=44826; B="KA01MU1745"; C="I20"; D="BODY SHOP"; E="WORK DONE DELIVERED"; F=15997; G="  INSURANCE"},
    @{A=44827; B="KA03MW6209"; C="BRIO"; D="PMS"; E="WORK DONE DELIVERED"; F=12586; G="CREDIT"},
    @{A=44827; B="KA51MG0886"; C="XCENT"; D="WIPER BLADE CHANGE"; E="WORK DONE DELIVERED"; F=739; G="CREDIT"},
    @{A=44827; B="KA51MD3348"; C="SUNNY"; D="RUNNING REPAIR"; E="WORK DONE DELIVERED"; F=6238; G=$null},
    @{A=44827; B="TS07F3598"; C="I20"; D="AC PROBLEM"; E="WORK DONE DELIVERED"; F=17808; G=$null},
    @{A=44827; B="KA04MM9589"; C="ECOSPORT"; D="POWER WINDOW SWITCH"; E="WORK DONE DELIVERED"; F=1376; G="CREDIT"},
    @{A=44827; B="KA01MU1745"; C="I20"; D="PMS"; E="WORK DONE DELIVERED"; F=8175; G=$null},
    @{A=44828; B="KA51MF8742"; C="ZEST"; D="PMS"; E="WORK DONE DELIVERED"; F=20257; G=$null},
    @{A=44828; B="KA10M2591"; C="POLO"; D="PMS"; E="WORK DONE DELIVERED"; F=25638; G=$null},
    @{A=44828; B="PB00RRC2164"; C="JEEP"; D="PMS"; E="WORK DONE DELIVERED"; F=5187; G=$null},
    @{A=44828; B="KA51N9050"; C="SWIFT"; D="PMS"; E="WORK DONE DELIVERED"; F=11317; G=$null},
    @{A=44828; B="KA03NL7739"; C="NANO"; D="BODY SHOP"; E="WORK DONE DELIVERED"; F=3052; G="CARD"},
    @{A=44828; B="KA05MV7398"; C="HEXA"; D="PMS"; E="WORK DONE DELIVERED"; F=8915; G=$null},
    @{A=44828; B="KA51ME7604"; C="BRIO"; D="PMS"; E="WORK DONE DELIVERED"; F=15012; G="P PAY"},
    @{A=44828; B="KA53P0246"; C="LINEA"; D="RUNNING REPAIR"; E="WORK DONE DELIVERED"; F=6500; G="GPAY"},
    @{A=44828; B="KA01MG7555"; C="BEAT"; D="PMS"; E="WORK DONE DELIVERED"; F=6112; G=$null},
    @{A=44828; B="KA51MD5564"; C="POLO"; D="BRAKE PAD CHANGE"; E="WORK DONE DELIVERED"; F=2309; G="CARD"}
)

$startRow = 423
$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    if ($row.G -ne $null) {
        $ws.Cells.Item($r, 7).Value = $row.G
    }
    $r = $r + 1
}

$ws.Range("G448").Select()